$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: pin 4 was disconnected (B14 cleared), a note explaining why is added in F14
$ws.Range("B14").Clear()
$ws.Range("F14").Value = "previously on PA4, but disconnected bc no need"

# Row 15: count initialization (pin 2) moved here
$ws.Range("B15").Value = 2

# Update the active selection to reflect where the user ended up (F15)
$ws.Range("F15").Select()
